$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C = "Done" header already present in C1; now fill in assignments.
# Row 7 is entered ahead of row 6 so that new shared-string entries land in
# the same order ("Grant" before "Steve") as the target workbook.
$assignments = @{
  2  = "X"
  3  = "Jordan"
  4  = "Ask Tas"
  5  = "Jordan"
  7  = "Grant"
  6  = "Steve"
  8  = "Jordan"
  9  = "Jordan"
  10 = "Steve"
  11 = "Justin"
  12 = "Steve"
  13 = "Grant"
  14 = "X"
  15 = "Steve"
  16 = "Justin"
  17 = "Justin"
  18 = "X"
  19 = "X"
  20 = "Grant"
  21 = "Grant"
  22 = "Grant"
  23 = "Justin"
}

$rowOrder = @(2,3,4,5,7,6,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)
foreach ($row in $rowOrder) {
  $ws.Cells.Item($row, 3).Value = $assignments[$row]
}

# Update view: scroll position and active selection.
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.ScrollColumn = 1
$ws.Range("C7").Select()
